$d = $word.ActiveDocument

# wdAlignParagraphJustify
$wdAlignParagraphJustify = 3

# --- Paragraph 1: "<tab><5 spaces>Defendant appeared in Court on {{ plea_trial_date }}..." ---
# Remove the leading tab run and the 5-space run entirely, then justify the paragraph.
$p1 = $d.Paragraphs.Item(12)
$p1Start = $p1.Range.Start
$leadRng = $d.Range($p1Start, $p1Start + 6)
if ($leadRng.Text -eq "`t     ") {
    $leadRng.Delete()
}
$p1.Alignment = $wdAlignParagraphJustify

# --- Paragraph 2: empty paragraph between the two blocks ---
# Just justify it.
$p2 = $d.Paragraphs.Item(13)
$p2.Alignment = $wdAlignParagraphJustify

# --- Paragraph 3: "<tab><5 spaces>The Court finds that the below-ordered conditions..." ---
# Keep the leading tab, remove only the 5 spaces that follow it, then justify the paragraph.
$p3 = $d.Paragraphs.Item(14)
$p3Start = $p3.Range.Start
$spacesRng = $d.Range($p3Start + 1, $p3Start + 6)
if ($spacesRng.Text -eq "     ") {
    $spacesRng.Delete()
}
$p3.Alignment = $wdAlignParagraphJustify
